$d = $word.ActiveDocument

# --- 1. "Programa Educativo: ____" line -> "Programa Educativo: {programa}" ---
# Locate the paragraph that contains the literal placeholder line so we can
# target just the underscore run (leave the "Programa Educativo: " label run
# untouched) and turn it into a single-underlined "{programa}" run.
$paras = $d.Paragraphs
$progParaIndex = 0
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "Programa Educativo:*____*") {
        $progParaIndex = $i
        break
    }
}

if ($progParaIndex -gt 0) {
    $label = "Programa Educativo: "
    $pRange = $paras.Item($progParaIndex).Range
    $underscoreStart = $pRange.Start + $label.Length
    $underscoreEnd = $pRange.End - 1  # exclude paragraph mark
    $underscoreRange = $d.Range($underscoreStart, $underscoreEnd)
    $underscoreRange.Font.Underline = 1
    $underscoreRange.Text = "{programa}"
}

# --- 2. Collapse the split "{ai-" / "N" / "}" runs into single "{ai-N}" runs ---
foreach ($n in 2,3,4,5,6) {
    $token = "{ai-$n}"
    $d.Content.Find.Execute($token, $true, $false, $false, $false, $false, $true, 1, $false, $token, 2) | Out-Null
}

Write-Output "done"
